$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 377
$ws1.Range("F6").Value = 661
$ws1.Range("F8").Value = 2091
$ws1.Range("F9").Value = 10847
$ws1.Range("F15").Value = 9062
$ws1.Range("F18").Value = 5321
$ws1.Range("F20").Value = 3371

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 377
$ws4.Range("F6").Value = 661
$ws4.Range("F9").Value = 2091
$ws4.Range("F12").Value = 10847
$ws4.Range("F18").Value = 9062
$ws4.Range("F21").Value = 5321
$ws4.Range("F23").Value = 3371
